$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2026-01-09T12:23:37.609996+00:00"
$ws.Range("B2").Value = "GET"
$ws.Range("C2").Value = "/api/inspiration-boards/products/6960ee6adfdda64b2c01ef32"
$ws.Range("D2").Value = "Inspiration Boards"
$ws.Range("E2").Value = 0.1571
$ws.Range("F2").Value = 200
$ws.Range("G2").Value = "6925b1e37b5978266363464e"
$ws.Range("H2").Value = ""

$ws.Range("A3").Value = "2026-01-09T12:23:37.759328+00:00"
$ws.Range("B3").Value = "GET"
$ws.Range("C3").Value = "/api/inspiration-boards/boards/6960eb294ad6a4df36746c2b"
$ws.Range("D3").Value = "Inspiration Boards"
$ws.Range("E3").Value = 0.0946
$ws.Range("F3").Value = 200
$ws.Range("G3").Value = "6925b1e37b5978266363464e"
$ws.Range("H3").Value = ""
